# add dcl revoke and reset password case
# (adds protocol_016 "desc" test case row to the Protocol test sheet)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 17 values -----------------------------------------------------
$ws.Range("A17").Value = "protocol_016"
$ws.Range("B17").Value = "y"
$ws.Range("C17").Value = "查看表的描述信息desc"
$ws.Range("D17").Value = "Protocol"
$ws.Range("E17").Value = "protocol_016"
$ws.Range("H17").Value = "create table protocol_016(id int auto_increment, acct_no varchar(20), tx_date date not null default '2020-10-01', tx_time time, tx_type int default 1, tx_status boolean, tx_amt double, tx_location varchar(255) default 'BJ', id_card_no varchar(18), phone varchar(11) not null, primary key(id, acct_no, tx_date))"
$ws.Range("I17").Value = "desc protocol_016"
$ws.Range("J17").Value = "src/test/resources/io.dingodb.test/testdata/mysqlcases/protocol/expectedresult/protocol_016.csv"
$ws.Range("K17").Value = "csv_containsAll"

# --- Formatting: mirror the style used by the rest of the table --------
$ws.Range("A17:B17").NumberFormat = "@"
$ws.Range("C17").NumberFormat = "@"
$ws.Range("D17:E17").NumberFormat = "@"
$ws.Range("H17").NumberFormat = "@"
$ws.Range("I17:J17").NumberFormat = "@"
$ws.Range("K17").NumberFormat = "@"

# --- Move the active selection like the source workbook (I21) ----------
$ws.Range("I21").Select()
